# Update Name of Algo
# Applies updated RandomForest imputation values to column B of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    12 = 4.9758
    32 = 6.551699999999998
    36 = 9.014900000000006
    38 = 5.299499999999999
    46 = 6.327400000000003
    54 = 4.641099999999999
    55 = 5.029499999999999
    67 = 5.925599999999995
    69 = 5.297099999999998
    72 = 5.129600000000003
    91 = 5.718400000000005
    99 = 4.582999999999998
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
